# Regenerate save_data column G ("K") values: use K (strikeouts) instead of
# Strike# for each outing row, per commit "regen save_data to use K instead
# of Strike#, regen std/mean, calc and write s_vals".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 2
    6  = 0
    7  = 0
    8  = 1
    9  = 2
    10 = 2
    11 = 1
    12 = 1
    13 = 0
    14 = 2
    15 = 2
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 1
    24 = 2
    25 = 2
    26 = 1
    27 = 0
    28 = 1
    29 = 0
    30 = 0
    31 = 2
    32 = 2
    33 = 2
    34 = 0
    35 = 1
    36 = 3
    37 = 1
    38 = 0
    39 = 0
    40 = 0
    41 = 1
    42 = 1
    43 = 1
    44 = 3
    45 = 0
    46 = 5
    47 = 1
    48 = 1
    49 = 2
    50 = 0
    51 = 3
    52 = 0
    53 = 2
    54 = 5
    55 = 1
    56 = 1
    57 = 1
    58 = 0
    59 = 3
    60 = 2
    61 = 1
    62 = 3
    63 = 1
    64 = 3
    65 = 2
    66 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
